# Apply the "Add new NCC with feature" update:
#  - ACR sheet (sheet2 / index 2): recompute all 101 columns x 5 rows.
#    Rows 1,2,3,5 split the genuine/impostor boundary between column 50 (AX) and 51 (AY).
#    Row 4 shifts that boundary one column earlier (between column 49 (AW) and 50 (AX)).
#  - FAR sheet (sheet3 / index 3): AX4 (row 4, col 50) flips from 1 to 0.
#  - FRR sheet (sheet4 / index 4): AX4 (row 4, col 50) flips from 0 to 1.

$wb = $excel.ActiveWorkbook

$lowValue  = 0.23255813953488372
$highValue = 0.76744186046511631

$totalCols = 101

# --- ACR sheet ---
$acr = $wb.Worksheets.Item("ACR")

for ($row = 1; $row -le 5; $row++) {
    if ($row -eq 4) {
        $splitCol = 49   # AW/AX boundary: AX4 already belongs to the "high" group
    } else {
        $splitCol = 50   # AX/AY boundary
    }

    for ($col = 1; $col -le $totalCols; $col++) {
        if ($col -le $splitCol) {
            $acr.Cells.Item($row, $col).Value = $lowValue
        } else {
            $acr.Cells.Item($row, $col).Value = $highValue
        }
    }
}

# --- FAR sheet ---
$far = $wb.Worksheets.Item("FAR")
$far.Cells.Item(4, 50).Value = 0

# --- FRR sheet ---
$frr = $wb.Worksheets.Item("FRR")
$frr.Cells.Item(4, 50).Value = 1

Write-Output "Applied NCC with feature EER train3test5 update"
